# Weekly price-log update: insert three new daily records (week of
# 2021-11-16, serial 44516) for "Feria Lagunitas de Puerto Montt - Lechuga"
# ahead of the existing history, pushing the prior rows (358:381) down to
# (361:384).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at 358, shifting everything currently at 358:381 down
# to 361:384 (formats/styles carry down with the shift, same as Excel UI).
$ws.Rows("358:360").Insert()

$newRows = @(
  @{ r = 358; D = 44516; H = "Escarola"; I = "Primera"; J = 200; K = 10000; L = 10000; M = 10000; N = "$/caja 15 unidades"; O = "Región de Coquimbo";   P = 667; Q = 15 },
  @{ r = 359; D = 44516; H = "Escarola"; I = "Segunda"; J = 200; K = 8000;  L = 8000;  M = 8000;  N = "$/caja 18 unidades"; O = "Región de Coquimbo";   P = 444; Q = 18 },
  @{ r = 360; D = 44516; H = "Marina";   I = "Primera"; J = 300; K = 8000;  L = 8500;  M = 8250;  N = "$/caja 15 unidades"; O = "Región Metropolitana"; P = 550; Q = 15 }
)

foreach ($row in $newRows) {
  $r = $row.r
  $ws.Range("A$r").Value = 4
  $ws.Range("B$r").Value = "Feria Lagunitas de Puerto Montt"
  $ws.Range("C$r").Value = "Los Lagos"
  $ws.Range("D$r").Value = $row.D
  $ws.Range("E$r").Value = 10
  $ws.Range("F$r").Value = 100112033
  $ws.Range("G$r").Value = "Lechuga"
  $ws.Range("H$r").Value = $row.H
  $ws.Range("I$r").Value = $row.I
  $ws.Range("J$r").Value = $row.J
  $ws.Range("K$r").Value = $row.K
  $ws.Range("L$r").Value = $row.L
  $ws.Range("M$r").Value = $row.M
  $ws.Range("N$r").Value = $row.N
  $ws.Range("O$r").Value = $row.O
  $ws.Range("P$r").Value = $row.P
  $ws.Range("Q$r").Value = $row.Q
  $ws.Range("R$r").Value = "Hortaliza"
}

Write-Output ("Updated dimension: " + $ws.UsedRange.Address())
